$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imports")

# "Fix incorrect external import": the ENVO row's "urban area" entry had
# the wrong id (ENVO:01000856 instead of ENVO:00000856) - correct it.
$ws.Range("D9").Value = "environment [ENVO:01000254]; rural area [ENVO:01000772]; urban area [ENVO:00000856]; research facility [ENVO:00000469]; park [ENVO:00000562]; forest [ENVO:00000111]; beach [ENVO:00000091]; grassland [ENVO:00000106]; road [ENVO:00000064]"

# Leave the sheet scrolled/selected near the corrected cell, matching the
# state the workbook was left in after making the fix.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E9").Select()
